$wb = $excel.ActiveWorkbook

# --- Work on the "Production" sheet (will become "Card Display") ---
$ws = $wb.Worksheets.Item("Production")

# The card layout used to have blocks of 9,8,9,9,8 columns (A:I, J:Q, R:Z, AA:AI, AJ:AQ).
# Change it to 5 even blocks of 9 columns each by inserting one column into the
# 2nd block (before the old column Q) and one column into the 5th/last block
# (before the old column AQ, which after the first insert sits at AR).
$ws.Columns("Q").Insert()
$ws.Columns("AR").Insert()

# Fix up the manual column page breaks, which used to sit right after each
# block (columns I, Q, Z, AI, AQ, plus a stray one at AZ). Remove the stale
# ones (their ids don't shift automatically on column insert) and re-add
# them after the new block boundaries (I, R, AA, AJ, AS).
$ws.Columns("J").PageBreak = -4142
$ws.Columns("R").PageBreak = -4142
$ws.Columns("AA").PageBreak = -4142
$ws.Columns("AJ").PageBreak = -4142
$ws.Columns("AR").PageBreak = -4142
$ws.Columns("BA").PageBreak = -4142

$ws.Columns("J").PageBreak = -4135
$ws.Columns("S").PageBreak = -4135
$ws.Columns("AB").PageBreak = -4135
$ws.Columns("AK").PageBreak = -4135
$ws.Columns("AT").PageBreak = -4135

# The selected cell was the top-left of the last card block; move it along
# with that block's new position.
$ws.Range("AK5").Select()

# --- Rename sheets and drop the unused "Operations" sheet ---
$wb.Worksheets.Item("Production").Name = "Card Display"
$wb.Worksheets.Item("Sorting").Name = "List Display"
$wb.Worksheets.Item("Operations").Delete()

# Keep "Card Display" as the active tab.
$wb.Worksheets.Item("Card Display").Activate()
